$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting existing C:I (Mutual Fund..QoQ) to D:J
$ws.Columns.Item(3).Insert()

# Copy header formatting from the (now shifted) Mutual Fund header cell to the new Industry header cell
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)
$ws.Cells.Item(1, 3).Value = "Industry"

# Populate Industry values for each data row
$industries = @{
    2 = "Pharmaceuticals & Biotechnology"
    3 = "Power"
    4 = "Construction"
    5 = "Chemicals & Petrochemicals"
    6 = "Petroleum Products"
    7 = "Metals & Minerals Trading"
    8 = "Auto Components"
    9 = "Agricultural Food & other Products"
    10 = "Construction"
    11 = "Entertainment"
    12 = "Finance"
    13 = "Power"
    14 = "Realty"
    15 = "Finance"
    16 = "Pharmaceuticals & Biotechnology"
    17 = "Finance"
    18 = "Textiles & Apparels"
    19 = "Fertilizers & Agrochemicals"
    20 = "Retailing"
    21 = "Automobiles"
    22 = "Financial Technology (Fintech)"
    23 = "Industrial Manufacturing"
    24 = "Telecom - Services"
    25 = "Transport Infrastructure"
    26 = "Food Products"
    27 = "Textiles & Apparels"
    28 = "Leisure Services"
    29 = "Industrial Products"
    30 = "Electrical Equipment"
    31 = "Aerospace & Defense"
    32 = "Gas"
    33 = "Construction"
    34 = "Minerals & Mining"
    35 = "Retailing"
    36 = "Healthcare Equipment & Supplies"
    37 = "Insurance"
    38 = "Industrial Products"
    39 = "Industrial Products"
    40 = "Pharmaceuticals & Biotechnology"
    41 = "Telecom - Services"
    42 = "Textiles & Apparels"
    43 = "Petroleum Products"
    44 = "Insurance"
    45 = "IT - Services"
    46 = "Transport Services"
    47 = "Leisure Services"
    48 = "Agricultural Food & other Products"
    49 = "Pharmaceuticals & Biotechnology"
    50 = "Realty"
    51 = "Food Products"
    52 = "Insurance"
    53 = "Fertilizers & Agrochemicals"
    54 = "IT - Software"
    55 = "Chemicals & Petrochemicals"
    56 = "Finance"
    57 = "Pharmaceuticals & Biotechnology"
    58 = "Leisure Services"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $industries[$row]
}
